# Snippet-extractor metadata workbook update
# - Fix mislabeled snippet mappings (RemoveDuplicateResult -> RemoveDuplicatesResult,
#   Image -> Shape for the saveAsPicture row)
# - Add two new rows to the "Snippets" table documenting the textbox / shape-text
#   snippet (excel-shape-textboxes): ShapeCollection.addTextBox / TextFrame.deleteText

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fixing mislabeled snippet mappings -----------------------------------
$ws.Range("A107").Value = "RemoveDuplicatesResult"
$ws.Range("A121").Value = "Shape"

# --- Snippet for textboxes and shapes --------------------------------------
$tbl = $ws.ListObjects.Item(1)

# Row 166: ShapeCollection.addTextBox
$tbl.ListRows.Add() | Out-Null
$ws.Range("A166").Value = "ShapeCollection"
$ws.Range("B166").Value = "addTextBox"
$ws.Range("C166").Value = "excel-shape-textboxes"
$ws.Range("D166").Value = "createTextbox"

# Row 167: TextFrame.deleteText
$tbl.ListRows.Add() | Out-Null
$ws.Range("D167").Value = "deleteText"
$ws.Range("B167").Value = "deleteText"
$ws.Range("C167").Value = "excel-shape-textboxes"
$ws.Range("A167").Value = "TextFrame"

# Leave the selection on the last edited cell, as in the authored workbook.
$ws.Range("B167").Select()
